# [ADD] New normalize way
# Update the "Valor" (B) column values on the active worksheet to use the
# new normalization factor (divide previous values by 1.922779637004405).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value  = 0.21323296341891762
$ws.Range("B4").Value  = 0.2285857367850797
$ws.Range("B5").Value  = 0.42646592683783524
$ws.Range("B6").Value  = 0.4571714735701594
$ws.Range("B7").Value  = 0.6396988902567529
$ws.Range("B8").Value  = 0.6857572103552391
$ws.Range("B9").Value  = 0.9143429471403188
$ws.Range("B10").Value = 1.1429286839253985
